$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E16:E44) previously listed periods in descending
# order (1903 down to 1611). The database was refreshed/re-sorted so the
# periods now appear in ascending order (1611 up to 1903).
$periods = @(
    "1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Two "Valor Mora" amounts (column F) were swapped between the first and the
# last data row as part of the database refresh.
$ws.Range("F16").Value = 30000
$ws.Range("F44").Value = 12000
